# Ranking atualizado trivia outubro
#
# A new trivia column (J, "TRIVIA 8 - Memento") has scores now, so the
# TOTAL column (B) formulas need to widen their SUM range from C:I to
# C:J for every team row, and their cached results recalculated.
# Also refresh the sheet's remembered cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2's TOTAL has its own standalone formula.
$ws.Range("B2").Formula = "=SUM(C2:J2)"

# Rows 3-14 share one formula group; re-applying the widened formula to
# the whole B3:B14 block keeps them on a single shared formula.
$ws.Range("B3:B14").Formula = "=SUM(C3:J3)"

# The workbook was last left with the bottom-right pane's selection on
# D10 instead of the previous K7.
$ws.Range("D10").Select()

$wb.Save()
